$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.187.17'
$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("D3").Value = '1.585.35'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("E4").Value = '  -0.06%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '211.43'
$ws.Range("E5").Value = '  +1.26%  '
$ws.Range("E6").Value = '  +0.74%  '
$ws.Range("E7").Value = '  -0.05%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.246'
$ws.Range("E8").Value = '  +0.44%  '
$ws.Range("E9").Value = '  -0.59%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '19.23'
$ws.Range("E10").Value = '  -1.69%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0849'
$ws.Range("E11").Value = '  +0.75%  '
$ws.Range("D12").Value = '1.808.80'
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("D13").Value = '1.588.34'
$ws.Range("E13").Value = '  +0.94%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '4.02'
$ws.Range("E14").Value = '  -0.91%  '
$ws.Range("E15").Value = '  +0.71%  '
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("D17").Value = '26.175.08'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("E18").Value = '  -0.28%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '7.34'
$ws.Range("E19").Value = '  +0.70%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '212.26'
$ws.Range("E20").Value = '  +1.67%  '
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("E22").Value = '  -0.60%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '2.17'
$ws.Range("E23").Value = '  +0.81%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '8.95'
$ws.Range("E24").Value = '  +1.40%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '143.43'
$ws.Range("E25").Value = '  -0.41%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("E28").Value = '  -0.63%  '
$ws.Range("E29").Value = '  -0.51%  '
$ws.Range("E30").Value = '  -2.06%  '
$ws.Range("E31").Value = '  +0.86%  '
$ws.Range("E32").Value = '  -0.76%  '
$ws.Range("D33").Value = '1.335.96'
$ws.Range("E33").Value = '  +4.30%  '
$ws.Range("E34").Value = '  -1.94%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  -1.17%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.580'
$ws.Range("E37").Value = '  -5.18%  '
$ws.Range("E38").Value = '  -0.44%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.817'
$ws.Range("E39").Value = '  +1.00%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '5.80'
$ws.Range("E40").Value = '  +3.92%  '
$ws.Range("E41").Value = '  -0.05%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.964'
$ws.Range("E42").Value = '  -12.95%  '
$ws.Range("E43").Value = '  +0.60%  '
$ws.Range("E44").Value = '  +0.42%  '
$ws.Range("D45").Value = '1.721.11'
$ws.Range("E45").Value = '  +0.55%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '60.91'
$ws.Range("E46").Value = '  -2.16%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '85.87'
$ws.Range("E47").Value = '  -3.14%  '
$ws.Range("E48").Value = '  -0.65%  '
$ws.Range("E49").Value = '  -1.80%  '
$ws.Range("E50").Value = '  -2.30%  '
$ws.Range("E51").Value = '  -0.76%  '
